# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header in G1 = "K") previously held the raw
# Strike# counts. This recalculates/overwrites the K column with the new
# s_vals that replace the old Strike# derived numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values for rows 2..66, in order.
$kValues = @(
    1,2,1,2,0,0,0,0,1,0,1,1,1,2,0,1,1,0,0,1,0,0,1,2,2,0,1,1,1,0,0,0,0,0,0,
    2,0,0,1,3,1,1,1,0,3,0,1,0,2,1,1,1,2,2,0,0,1,0,1,1,1,3,1,1,2
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Count; $i++) {
    $row = $startRow + $i
    $ws.Range("G$row").Value = $kValues[$i]
}
